$d = $word.ActiveDocument

# The paragraph currently holds the literal text "<id>p029r_1</id>" split
# across three runs: "<id>" (Courier New / color 7f6000 / sz 18), "p029r_1"
# (default formatting), "</id>" (Courier New / color 7f6000 / sz 18).
# Collapse them into a single run containing the full text, adopting the
# formatting of the run the replacement starts in (the "<id>" run), by doing
# a Find/Replace-in-place of the exact same text.
$d.Content.Find.Execute(
    "<id>p029r_1</id>",  # FindText
    $false,              # MatchCase
    $false,              # MatchWholeWord
    $false,              # MatchWildcards
    $false,              # MatchSoundsLike
    $false,              # MatchAllWordForms
    $true,               # Forward
    1,                   # Wrap (wdFindContinue)
    $false,              # Format
    "<id>p029r_1</id>",  # ReplaceWith
    2                    # Replace (wdReplaceAll)
) | Out-Null
